# fm12 data correction: add LocMinDed6All column (BQ) and update CondNumber (W) for rows 15-32
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column BQ (69)
$ws.Cells.Item(1, 69).Value = "LocMinDed6All"

# New LocMinDed6All values for rows 2-63 (column BQ / 69)
$bq = @{
    2 = 0
    3 = 0
    4 = 0
    5 = 0
    6 = 0
    7 = 0
    8 = 0
    9 = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 25000
    23 = 25000
    24 = 25000
    25 = 78783.03
    26 = 25000
    27 = 89642.04
    28 = 25000
    29 = 85681.983
    30 = 140183.82
    31 = 25000
    32 = 462684.45
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 0
    39 = 0
    40 = 0
    41 = 0
    42 = 0
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 0
    49 = 0
    50 = 0
    51 = 0
    52 = 0
    53 = 0
    54 = 0
    55 = 0
    56 = 0
    57 = 0
    58 = 0
    59 = 0
    60 = 0
    61 = 0
    62 = 0
    63 = 0
}

foreach ($r in $bq.Keys) {
    $ws.Cells.Item($r, 69).Value = $bq[$r]
}

# Updated CondNumber values (column W / 23) for rows 15-32
$wvals = @{
    15 = 1160
    16 = 1160
    17 = 1160
    18 = 1160
    19 = 1160
    20 = 1160
    21 = 1160
    22 = 1161
    23 = 1161
    24 = 1161
    25 = 1161
    26 = 1161
    27 = 1161
    28 = 1161
    29 = 1161
    30 = 1161
    31 = 1161
    32 = 1163
}

foreach ($r in $wvals.Keys) {
    $ws.Cells.Item($r, 23).Value = $wvals[$r]
}

# Update the active selection to match the saved view state
[void]$ws.Range("AH15").Select()
